$wb = $excel.ActiveWorkbook

# --- Add the new "metadata" worksheet, positioned right after "data" ---
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "metadata"
$wsNew.Move($null, $wb.Worksheets.Item("data"))

# Move() invalidates the earlier index-based worksheet reference, so
# re-fetch both sheets by name before doing any further work.
$ws1 = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("metadata")

# --- Populate "metadata" header row ---
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# --- Populate "metadata" data row ---
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Congenital hyperinsulinism"
$ws2.Range("C2").Value = 308

# "data_version" (2.5) is a text value in the source data, not a number -
# force text storage so it round-trips as a string rather than being
# auto-coerced to numeric.
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "2.5"

$ws2.Range("E2").Value = "2021-01-29T10:25:58.504503Z"
$ws2.Range("F2").Value = "2021-10-05 14:19:43.283292"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/308/?format=json"

# --- Reuse the "data" sheet's header style (bold/bordered/centered) for the
# new "metadata" header row and its A2 index cell, instead of minting a
# fresh style definition. Copy/PasteSpecial(formats) after the values are
# already in place. ---
$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# --- Refresh the "data" sheet's time_taken column (F2:F22) timestamps ---
$newTimes = @(
  "2021-10-05 14:19:43.287191",
  "2021-10-05 14:19:43.287200",
  "2021-10-05 14:19:43.287203",
  "2021-10-05 14:19:43.287206",
  "2021-10-05 14:19:43.287209",
  "2021-10-05 14:19:43.287211",
  "2021-10-05 14:19:43.287214",
  "2021-10-05 14:19:43.287216",
  "2021-10-05 14:19:43.287219",
  "2021-10-05 14:19:43.287222",
  "2021-10-05 14:19:43.287224",
  "2021-10-05 14:19:43.287227",
  "2021-10-05 14:19:43.287229",
  "2021-10-05 14:19:43.287232",
  "2021-10-05 14:19:43.287234",
  "2021-10-05 14:19:43.287237",
  "2021-10-05 14:19:43.287240",
  "2021-10-05 14:19:43.287242",
  "2021-10-05 14:19:43.287245",
  "2021-10-05 14:19:43.287248",
  "2021-10-05 14:19:43.287250"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
  $row = $i + 2
  $ws1.Cells.Item($row, 6).Value = $newTimes[$i]
}
